$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a number/percent
# (e.g. "0.61%") need NumberFormat forced to Text first, then restored to the
# workbook default style afterwards so no stray style index is left on the cell.

$ws.Range("B2").Value = "₹ 15,142"
$ws.Range("C2").Value = "₹ 74,831"
$ws.Range("E2").Value = "₹ 1,02,473"
$ws.Range("F2").Value = "₹ 24,84,858"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0.61%"
$ws.Range("G2").Style = "Normal"
$ws.Range("B3").Value = "₹ 48,238"
$ws.Range("C3").Value = "₹ 2,21,680"
$ws.Range("E3").Value = "₹ 3,07,418"
$ws.Range("F3").Value = "₹ 24,36,620"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "2.54%"
$ws.Range("G3").Style = "Normal"
$ws.Range("B4").Value = "₹ 52,763"
$ws.Range("C4").Value = "₹ 2,17,155"
$ws.Range("E4").Value = "₹ 3,07,418"
$ws.Range("F4").Value = "₹ 23,83,857"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "4.65%"
$ws.Range("G4").Style = "Normal"
$ws.Range("B5").Value = "₹ 57,713"
$ws.Range("C5").Value = "₹ 2,12,205"
$ws.Range("E5").Value = "₹ 3,07,418"
$ws.Range("F5").Value = "₹ 23,26,144"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "6.95%"
$ws.Range("G5").Style = "Normal"
$ws.Range("B6").Value = "₹ 63,127"
$ws.Range("C6").Value = "₹ 2,06,791"
$ws.Range("E6").Value = "₹ 3,07,418"
$ws.Range("F6").Value = "₹ 22,63,017"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "9.48%"
$ws.Range("G6").Style = "Normal"
$ws.Range("B7").Value = "₹ 69,048"
$ws.Range("C7").Value = "₹ 2,00,869"
$ws.Range("E7").Value = "₹ 3,07,418"
$ws.Range("F7").Value = "₹ 21,93,969"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "12.24%"
$ws.Range("G7").Style = "Normal"
$ws.Range("B8").Value = "₹ 75,526"
$ws.Range("C8").Value = "₹ 1,94,392"
$ws.Range("E8").Value = "₹ 3,07,418"
$ws.Range("F8").Value = "₹ 21,18,443"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "15.26%"
$ws.Range("G8").Style = "Normal"
$ws.Range("B9").Value = "₹ 82,610"
$ws.Range("C9").Value = "₹ 1,87,307"
$ws.Range("E9").Value = "₹ 3,07,418"
$ws.Range("F9").Value = "₹ 20,35,833"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "18.57%"
$ws.Range("G9").Style = "Normal"
$ws.Range("B10").Value = "₹ 90,360"
$ws.Range("C10").Value = "₹ 1,79,558"
$ws.Range("E10").Value = "₹ 3,07,418"
$ws.Range("F10").Value = "₹ 19,45,473"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "22.18%"
$ws.Range("G10").Style = "Normal"
$ws.Range("B11").Value = "₹ 98,836"
$ws.Range("C11").Value = "₹ 1,71,082"
$ws.Range("E11").Value = "₹ 3,07,418"
$ws.Range("F11").Value = "₹ 18,46,637"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "26.13%"
$ws.Range("G11").Style = "Normal"
$ws.Range("B12").Value = "₹ 1,08,108"
$ws.Range("C12").Value = "₹ 1,61,810"
$ws.Range("E12").Value = "₹ 3,07,418"
$ws.Range("F12").Value = "₹ 17,38,529"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "30.46%"
$ws.Range("G12").Style = "Normal"
$ws.Range("B13").Value = "₹ 1,18,249"
$ws.Range("C13").Value = "₹ 1,51,669"
$ws.Range("E13").Value = "₹ 3,07,418"
$ws.Range("F13").Value = "₹ 16,20,280"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "35.19%"
$ws.Range("G13").Style = "Normal"
$ws.Range("B14").Value = "₹ 1,29,342"
$ws.Range("C14").Value = "₹ 1,40,576"
$ws.Range("E14").Value = "₹ 3,07,418"
$ws.Range("F14").Value = "₹ 14,90,939"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "40.36%"
$ws.Range("G14").Style = "Normal"
$ws.Range("B15").Value = "₹ 1,41,475"
$ws.Range("C15").Value = "₹ 1,28,443"
$ws.Range("E15").Value = "₹ 3,07,418"
$ws.Range("F15").Value = "₹ 13,49,464"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "46.02%"
$ws.Range("G15").Style = "Normal"
$ws.Range("B16").Value = "₹ 1,54,746"
$ws.Range("C16").Value = "₹ 1,15,172"
$ws.Range("E16").Value = "₹ 3,07,418"
$ws.Range("F16").Value = "₹ 11,94,718"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "52.21%"
$ws.Range("G16").Style = "Normal"
$ws.Range("B17").Value = "₹ 1,69,262"
$ws.Range("C17").Value = "₹ 1,00,656"
$ws.Range("E17").Value = "₹ 3,07,418"
$ws.Range("F17").Value = "₹ 10,25,456"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "58.98%"
$ws.Range("G17").Style = "Normal"
$ws.Range("B18").Value = "₹ 1,85,140"
$ws.Range("C18").Value = "₹ 84,778"
$ws.Range("E18").Value = "₹ 3,07,418"
$ws.Range("F18").Value = "₹ 8,40,315"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "66.39%"
$ws.Range("G18").Style = "Normal"
$ws.Range("B19").Value = "₹ 2,02,508"
$ws.Range("C19").Value = "₹ 67,410"
$ws.Range("E19").Value = "₹ 3,07,418"
$ws.Range("F19").Value = "₹ 6,37,808"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "74.49%"
$ws.Range("G19").Style = "Normal"
$ws.Range("B20").Value = "₹ 2,21,504"
$ws.Range("C20").Value = "₹ 48,414"
$ws.Range("E20").Value = "₹ 3,07,418"
$ws.Range("F20").Value = "₹ 4,16,304"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "83.35%"
$ws.Range("G20").Style = "Normal"
$ws.Range("B21").Value = "₹ 2,42,283"
$ws.Range("C21").Value = "₹ 27,635"
$ws.Range("E21").Value = "₹ 3,07,418"
$ws.Range("F21").Value = "₹ 1,74,021"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "93.04%"
$ws.Range("G21").Style = "Normal"
$ws.Range("B22").Value = "₹ 1,74,021"
$ws.Range("C22").Value = "₹ 5,924"
$ws.Range("E22").Value = "₹ 2,04,945"
